$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.480.61"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "3.009.71"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'509.97"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'139.54"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.438"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'7.53"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +3.66%  "
$ws.Range("D12").Value = "3.520.85"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "'26.40"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("E15").Value = "  +6.78%  "
$ws.Range("D16").Value = "57.440.30"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "'6.21"
$ws.Range("E17").Value = "  +6.10%  "
$ws.Range("D18").Value = "3.005.41"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").Value = "'7.98"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "'331.20"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'0.499"
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("D24").Value = "'64.46"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "0.0₃0923"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "'6.82"
$ws.Range("E28").Value = "  +4.31%  "
$ws.Range("D29").Value = "'7.44"
$ws.Range("E29").Value = "  +7.02%  "
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").Value = "'1.20"
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("D32").Value = "'20.67"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'4.74"
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("D34").Value = "'154.25"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "'5.89"
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'24.49"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "'0.0683"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").Value = "3.039.17"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'37.32"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "'3.85"
$ws.Range("E41").Value = "  +6.74%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "2.285.84"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'0.990"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'6.03"
$ws.Range("E47").Value = "  +4.50%  "
$ws.Range("D48").Value = "'0.0240"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'19.49"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'1.86"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'0.0894"
$ws.Range("E51").Value = "  +2.20%  "
